# Auto-generated edit script: applies scheduled-runner market-data refresh
# to the Leviathan_Profits workbook (per-sheet leve profit calculations).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$updates = @(
    @("H51", 5558576),
    @("I51", 2830.2),
    @("K51", 2830.2),
    @("M51", -2346.2),
    @("H70", 4539.75),
    @("J70", 6215.4),
    @("L70", 18646.2),
    @("N70", -19186.2),
    @("H73", 4539.75),
    @("J73", 6215.4),
    @("L73", 18646.2),
    @("N73", -20518.2),
    @("H74", 5142.857),
    @("I74", 0),
    @("K74", 0),
    @("H77", 5142.857),
    @("I77", 0),
    @("K77", 0),
    @("H80", 1009.2727),
    @("J80", 1025.625),
    @("L80", 3076.875),
    @("N80", -5072.875),
    @("H83", 1009.2727),
    @("J83", 1025.625),
    @("L83", 9230.625),
    @("N83", -19214.625),
    @("H88", 92469.45),
    @("J88", 92469.45),
    @("L88", 92469.45),
    @("N88", -93281.45),
    @("H91", 92469.45),
    @("J91", 92469.45),
    @("L91", 92469.45),
    @("N91", -95277.45),
    @("H100", 3278.4707),
    @("I100", 2721.818),
    @("K100", 2721.818),
    @("M100", -2180.818),
    @("H104", 647.75),
    @("I104", 647.75),
    @("J104", 0),
    @("K104", 1943.25),
    @("L104", 0),
    @("M104", -196.25),
    @("H137", 57754.11),
    @("I137", 1962.1428),
    @("K137", 5886.428400000001),
    @("M137", -3336.428400000001),
    @("H138", 4169.04),
    @("J138", 4912.3125),
    @("L138", 14736.9375),
    @("N138", -25016.9375),
    @("H141", 45387.24),
    @("I141", 52607.11),
    @("K141", 157821.33),
    @("M141", -152641.33)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
$clears = @("M74", "M77", "N104")
foreach ($c in $clears) {
    $ws.Range($c).ClearContents()
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$updates = @(
    @("H32", 131081.03),
    @("I32", 121298.74),
    @("K32", 121298.74),
    @("M32", -121011.74),
    @("H132", 2597.04),
    @("I132", 2477.7144),
    @("K132", 7433.1432),
    @("M132", -4903.1432)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$updates = @(
    @("H22", 466),
    @("I22", 466),
    @("J22", 0),
    @("K22", 466),
    @("L22", 0),
    @("M22", -293)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
$clears = @("N22")
foreach ($c in $clears) {
    $ws.Range($c).ClearContents()
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$updates = @(
    @("H16", 1947.8889),
    @("I16", 1893.6923),
    @("J16", 2088.8),
    @("K16", 1893.6923),
    @("L16", 2088.8),
    @("M16", -1606.6923),
    @("N16", -2662.8),
    @("H111", 77925.336),
    @("J111", 77925.336),
    @("L111", 77925.336),
    @("N111", -86105.336),
    @("H113", 1947.8889),
    @("I113", 1893.6923),
    @("J113", 2088.8),
    @("K113", 1893.6923),
    @("L113", 2088.8),
    @("M113", 276.3077000000001),
    @("N113", -6428.8),
    @("H122", 1705),
    @("I122", 1306),
    @("K122", 3918),
    @("M122", -1468)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$updates = @(
    @("H2", 77),
    @("I2", 72.181816),
    @("J2", 94.666664),
    @("K2", 433.090896),
    @("L2", 567.999984),
    @("M2", -320.090896),
    @("N2", -793.999984),
    @("H74", 0),
    @("J74", 0),
    @("L74", 0),
    @("H77", 0),
    @("J77", 0),
    @("L77", 0)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
$clears = @("N74", "N77")
foreach ($c in $clears) {
    $ws.Range($c).ClearContents()
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$updates = @(
    @("H70", 6585.857),
    @("I70", 6390.364),
    @("K70", 6390.364),
    @("M70", -6120.364),
    @("H73", 6585.857),
    @("I73", 6390.364),
    @("K73", 6390.364),
    @("M73", -5454.364),
    @("H102", 4204.35),
    @("I102", 3663),
    @("K102", 3663),
    @("M102", -2041),
    @("H122", 2708.6428),
    @("I122", 2669.7273),
    @("J122", 2851.3333),
    @("K122", 8009.1819),
    @("L122", 8553.999899999999),
    @("M122", -5559.1819),
    @("N122", -13453.9999),
    @("H126", 3421),
    @("I126", 3128),
    @("K126", 9384),
    @("M126", -6914)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$updates = @(
    @("H7", 47428.57),
    @("I7", 64100),
    @("K7", 64100),
    @("M7", -63988),
    @("H93", 34841.4),
    @("I93", 1178.2),
    @("J93", 68504.60000000001),
    @("K93", 1178.2),
    @("L93", 68504.60000000001),
    @("M93", 69.79999999999995),
    @("N93", -71000.60000000001),
    @("H122", 11045.571),
    @("I122", 12103.167),
    @("J122", 4700),
    @("K122", 36309.501),
    @("L122", 14100),
    @("M122", -33859.501),
    @("N122", -19000),
    @("H126", 47428.57),
    @("I126", 64100),
    @("K126", 192300),
    @("M126", -189830),
    @("H136", 5576.846),
    @("I136", 4249.1665),
    @("J136", 6714.857),
    @("K136", 12747.4995),
    @("L136", 20144.571),
    @("M136", -10197.4995),
    @("N136", -25244.571)
)
foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
